# Apply the "Anonimyzed fedcore" update:
#  - rename the "fedcore" approach label to "approach" on both sheets
#  - add a thin top/bottom box border around the merged header cells:
#      the inner cell of each merged group keeps only top+bottom edges,
#      the last (rightmost) cell of the group additionally keeps the
#      right edge
#  - clear the stray empty cell G5 on the computational_comparison sheet

$wb = $excel.ActiveWorkbook

$xlContinuous = 1
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlPasteFormats = -4122

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two reusable border styles once (on quality_comparison!C1/D1) ---
# "inner" style: thin top + thin bottom, no left/right
$innerTemplate = $ws1.Range("C1")
$innerTemplate.ClearFormats()
$innerTemplate.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$innerTemplate.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# "last" style: same as inner, plus a thin right edge
$innerTemplate.Copy()
$lastTemplate = $ws1.Range("D1")
$lastTemplate.PasteSpecial($xlPasteFormats)
$lastTemplate.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# --- Reuse those two styles (via copy/paste of formats) on every other header cell ---
$innerTemplate.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$lastTemplate.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)

$innerTemplate.Copy()
$ws2.Range("F1").PasteSpecial($xlPasteFormats)
$lastTemplate.Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# --- Rename the "fedcore" column header to "approach" (anonymized) ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty cell ---
$ws2.Range("G5").ClearContents()

Write-Host "edit applied"
